$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 9895.17
$ws.Range("B6").Value = 9971.9500000000007
$ws.Range("C6").Value = 79.650000000000006
$ws.Range("D6").Value = 79.040000000000006
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -0.77
$ws.Range("G6").Value = 42612.674247685187
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"
$ws.Range("H6").Value = $false
